# feat: add 2022-Q4 data
#
# The existing "2021-Q3" sheet is renamed to "2022-Q4" and its data is
# replaced with the new quarter's figures. A brand-new sheet named
# "2021-Q3" is inserted right after it, carrying the data that used to
# live on the sheet we just renamed. The "总计" (totals) sheet gets a
# new row for the old quarter and its existing row is updated to point
# at the new quarter.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)
$q3sheet = $wb.Worksheets.Item(2)

# --- capture the data currently on the "2021-Q3" sheet before we overwrite it ---
$oldHeaders = @("基金代码", "基金名称", "基金金额", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$oldRow2 = @("012348", "天弘恒生科技指数型发起式证券投资基金（QDII）A", "12.21", "91.95", "4.36", "0.5324", 9)
$oldRow3 = @("012349", "天弘恒生科技指数型发起式证券投资基金（QDII）C", "5.09", "91.95", "4.36", "0.2219", 9)

# --- rename the current "2021-Q3" sheet to "2022-Q4" and add a fresh sheet "2021-Q3" after it ---
$q3sheet.Name = "2022-Q4"
$q4NewOld = $wb.Worksheets.Add($null, $q3sheet)
$q4NewOld.Name = "2021-Q3"

# --- fill the new "2021-Q3" sheet with the data that used to be on that tab ---
$ws3 = $q4NewOld
$ws3.Range("B1").Value = $oldHeaders[0]
$ws3.Range("C1").Value = $oldHeaders[1]
$ws3.Range("D1").Value = $oldHeaders[2]
$ws3.Range("E1").Value = $oldHeaders[3]
$ws3.Range("F1").Value = $oldHeaders[4]
$ws3.Range("G1").Value = $oldHeaders[5]
$ws3.Range("H1").Value = $oldHeaders[6]
$ws3.Range("B1:H1").Style = $total.Range("B1").Style

$ws3.Range("A2").Value = 0
$ws3.Range("A2").Style = $total.Range("B1").Style
$ws3.Range("B2").NumberFormat = "@"
$ws3.Range("B2").Value = $oldRow2[0]
$ws3.Range("C2").NumberFormat = "@"
$ws3.Range("C2").Value = $oldRow2[1]
$ws3.Range("D2").NumberFormat = "@"
$ws3.Range("D2").Value = $oldRow2[2]
$ws3.Range("E2").NumberFormat = "@"
$ws3.Range("E2").Value = $oldRow2[3]
$ws3.Range("F2").NumberFormat = "@"
$ws3.Range("F2").Value = $oldRow2[4]
$ws3.Range("G2").NumberFormat = "@"
$ws3.Range("G2").Value = $oldRow2[5]
$ws3.Range("H2").Value = $oldRow2[6]

$ws3.Range("A3").Value = 1
$ws3.Range("A3").Style = $total.Range("B1").Style
$ws3.Range("B3").NumberFormat = "@"
$ws3.Range("B3").Value = $oldRow3[0]
$ws3.Range("C3").NumberFormat = "@"
$ws3.Range("C3").Value = $oldRow3[1]
$ws3.Range("D3").NumberFormat = "@"
$ws3.Range("D3").Value = $oldRow3[2]
$ws3.Range("E3").NumberFormat = "@"
$ws3.Range("E3").Value = $oldRow3[3]
$ws3.Range("F3").NumberFormat = "@"
$ws3.Range("F3").Value = $oldRow3[4]
$ws3.Range("G3").NumberFormat = "@"
$ws3.Range("G3").Value = $oldRow3[5]
$ws3.Range("H3").Value = $oldRow3[6]

# --- overwrite the renamed sheet ("2022-Q4") with the new quarter's data ---
$ws2 = $wb.Worksheets.Item("2022-Q4")
$ws2.Range("B1").Value = "基金代码"
$ws2.Range("C1").Value = "基金名称"
$ws2.Range("D1").Value = "基金规模"
$ws2.Range("E1").Value = "股票总仓位"
$ws2.Range("F1").Value = "仓位占比"
$ws2.Range("G1").Value = "持有市值(亿元)"
$ws2.Range("H1").Value = "仓位排名"

$ws2.Range("A2").Value = 0
$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "012348"
$ws2.Range("C2").NumberFormat = "@"
$ws2.Range("C2").Value = "天弘恒生科技指数（QDII）A"
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "39.65"
$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "93.67"
$ws2.Range("F2").NumberFormat = "@"
$ws2.Range("F2").Value = "4.75"
$ws2.Range("G2").NumberFormat = "@"
$ws2.Range("G2").Value = "1.8834"
$ws2.Range("H2").Value = 8

$ws2.Range("A3").Value = 1
$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "012349"
$ws2.Range("C3").NumberFormat = "@"
$ws2.Range("C3").Value = "天弘恒生科技指数（QDII）C"
$ws2.Range("D3").NumberFormat = "@"
$ws2.Range("D3").Value = "37.52"
$ws2.Range("E3").NumberFormat = "@"
$ws2.Range("E3").Value = "93.67"
$ws2.Range("F3").NumberFormat = "@"
$ws2.Range("F3").Value = "4.75"
$ws2.Range("G3").NumberFormat = "@"
$ws2.Range("G3").Value = "1.7822"
$ws2.Range("H3").Value = 8

# --- update the "总计" (totals) sheet: row 2 now reflects 2022-Q4, row 3 keeps 2021-Q3 ---
$total.Range("B2").Value = "2022-Q4"
$total.Range("D2").Value = 3.67

$total.Range("A3").Value = 1
$total.Range("A3").Style = $total.Range("A2").Style
$total.Range("B3").Value = "2021-Q3"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.75
